$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row containing "Christian Watson" in column A and delete the entire row,
# shifting subsequent rows up.
$found = $ws.Range("A:A").Find("Christian Watson")
if ($found -ne $null) {
    $found.EntireRow.Delete()
}

# Update the view state to match the post-edit selection/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("G49").Select()
